# Apply the new table style ({0903D687-0B98-4071-9048-55D6C756C9BF}) to the
# three tables that previously used the deck's custom table style
# ({CAC69970-249F-4958-9B4B-790541CFD7EB}). These live on slides 14, 15 and 16,
# each as the sole shape (a single graphicFrame containing the table).

$p = $ppt.ActivePresentation

$oldStyleId = "{CAC69970-249F-4958-9B4B-790541CFD7EB}"
$newStyleId = "{0903D687-0B98-4071-9048-55D6C756C9BF}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
